$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 23.78418466666666
$ws.Cells.Item(2, 8).Value = 71.352554
$ws.Cells.Item(2, 9).Value = 0.06460357633592957
$ws.Cells.Item(2, 10).Value = 0.06460357633592959
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.06778666666666666
$ws.Cells.Item(2, 14).Value = 0.20336
$ws.Cells.Item(2, 15).Value = 0.01026668284214455
$ws.Cells.Item(2, 16).Value = 0.01026668284214455
$ws.Cells.Item(2, 17).Value = 1.612250597937777
$ws.Cells.Item(2, 18).Value = 14.51025538144
$ws.Cells.Item(2, 19).Value = 0.0006632644287092638
$ws.Cells.Item(2, 20).Value = 0.000663264428709264
$ws.Cells.Item(3, 7).Value = 23.78418466666666
$ws.Cells.Item(3, 8).Value = 71.352554
$ws.Cells.Item(3, 9).Value = 0.06460357633592957
$ws.Cells.Item(3, 10).Value = 0.06460357633592959
$ws.Cells.Item(3, 15).Value = 0.01567037284022157
$ws.Cells.Item(3, 16).Value = 0.01567037284022157
$ws.Cells.Item(3, 17).Value = 2.460830666536666
$ws.Cells.Item(3, 18).Value = 22.14747599883
$ws.Cells.Item(3, 19).Value = 0.001012362127995732
$ws.Cells.Item(3, 20).Value = 0.001012362127995732
$ws.Cells.Item(4, 7).Value = 23.78418466666666
$ws.Cells.Item(4, 8).Value = 71.352554
$ws.Cells.Item(4, 9).Value = 0.06460357633592957
$ws.Cells.Item(4, 10).Value = 0.06460357633592959
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.08443199999999999
$ws.Cells.Item(4, 14).Value = 0.253296
$ws.Cells.Item(4, 15).Value = 0.01278771487600239
$ws.Cells.Item(4, 16).Value = 0.01278771487600239
$ws.Cells.Item(4, 17).Value = 2.008146279776
$ws.Cells.Item(4, 18).Value = 18.073316517984
$ws.Cells.Item(4, 19).Value = 0.0008261321141539226
$ws.Cells.Item(4, 20).Value = 0.0008261321141539227
$ws.Cells.Item(5, 7).Value = 23.78418466666666
$ws.Cells.Item(5, 8).Value = 71.352554
$ws.Cells.Item(5, 9).Value = 0.06460357633592957
$ws.Cells.Item(5, 10).Value = 0.06460357633592959
$ws.Cells.Item(5, 13).Value = 6.346903333333334
$ws.Cells.Item(5, 14).Value = 19.04071
$ws.Cells.Item(5, 15).Value = 0.9612752294416316
$ws.Cells.Item(5, 16).Value = 0.9612752294416316
$ws.Cells.Item(5, 17).Value = 150.9559209414822
$ws.Cells.Item(5, 18).Value = 1358.60328847334
$ws.Cells.Item(5, 19).Value = 0.06210181766507066
$ws.Cells.Item(5, 20).Value = 0.06210181766507068
$ws.Cells.Item(6, 9).Value = 0.3773880863345054
$ws.Cells.Item(6, 10).Value = 0.3773880863345054
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.06778666666666666
$ws.Cells.Item(6, 14).Value = 0.20336
$ws.Cells.Item(6, 15).Value = 0.01026668284214455
$ws.Cells.Item(6, 16).Value = 0.01026668284214455
$ws.Cells.Item(6, 17).Value = 9.418118970435556
$ws.Cells.Item(6, 18).Value = 84.76307073392
$ws.Cells.Item(6, 19).Value = 0.003874523790800233
$ws.Cells.Item(6, 20).Value = 0.003874523790800233
$ws.Cells.Item(7, 9).Value = 0.3773880863345054
$ws.Cells.Item(7, 10).Value = 0.3773880863345054
$ws.Cells.Item(7, 15).Value = 0.01567037284022157
$ws.Cells.Item(7, 16).Value = 0.01567037284022157
$ws.Cells.Item(7, 19).Value = 0.005913812018319425
$ws.Cells.Item(7, 20).Value = 0.005913812018319425
$ws.Cells.Item(8, 9).Value = 0.3773880863345054
$ws.Cells.Item(8, 10).Value = 0.3773880863345054
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.08443199999999999
$ws.Cells.Item(8, 14).Value = 0.253296
$ws.Cells.Item(8, 15).Value = 0.01278771487600239
$ws.Cells.Item(8, 16).Value = 0.01278771487600239
$ws.Cells.Item(8, 17).Value = 11.730782173168
$ws.Cells.Item(8, 18).Value = 105.577039558512
$ws.Cells.Item(8, 19).Value = 0.004825931245645829
$ws.Cells.Item(8, 20).Value = 0.004825931245645829
$ws.Cells.Item(9, 9).Value = 0.3773880863345054
$ws.Cells.Item(9, 10).Value = 0.3773880863345054
$ws.Cells.Item(9, 13).Value = 6.346903333333334
$ws.Cells.Item(9, 14).Value = 19.04071
$ws.Cells.Item(9, 15).Value = 0.9612752294416316
$ws.Cells.Item(9, 16).Value = 0.9612752294416316
$ws.Cells.Item(9, 17).Value = 881.8237217818745
$ws.Cells.Item(9, 18).Value = 7936.413496036871
$ws.Cells.Item(9, 19).Value = 0.36277381927974
$ws.Cells.Item(9, 20).Value = 0.36277381927974
$ws.Cells.Item(10, 7).Value = 150.629115
$ws.Cells.Item(10, 8).Value = 451.887345
$ws.Cells.Item(10, 9).Value = 0.4091449703110563
$ws.Cells.Item(10, 10).Value = 0.4091449703110563
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.06778666666666666
$ws.Cells.Item(10, 14).Value = 0.20336
$ws.Cells.Item(10, 15).Value = 0.01026668284214455
$ws.Cells.Item(10, 16).Value = 0.01026668284214455
$ws.Cells.Item(10, 17).Value = 10.2106456088
$ws.Cells.Item(10, 18).Value = 91.89581047919999
$ws.Cells.Item(10, 19).Value = 0.004200561646642264
$ws.Cells.Item(10, 20).Value = 0.004200561646642264
$ws.Cells.Item(11, 7).Value = 150.629115
$ws.Cells.Item(11, 8).Value = 451.887345
$ws.Cells.Item(11, 9).Value = 0.4091449703110563
$ws.Cells.Item(11, 10).Value = 0.4091449703110563
$ws.Cells.Item(11, 15).Value = 0.01567037284022157
$ws.Cells.Item(11, 16).Value = 0.01567037284022157
$ws.Cells.Item(11, 17).Value = 15.584841383475
$ws.Cells.Item(11, 18).Value = 140.263572451275
$ws.Cells.Item(11, 19).Value = 0.006411454230475637
$ws.Cells.Item(11, 20).Value = 0.006411454230475637
$ws.Cells.Item(12, 7).Value = 150.629115
$ws.Cells.Item(12, 8).Value = 451.887345
$ws.Cells.Item(12, 9).Value = 0.4091449703110563
$ws.Cells.Item(12, 10).Value = 0.4091449703110563
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.08443199999999999
$ws.Cells.Item(12, 14).Value = 0.253296
$ws.Cells.Item(12, 15).Value = 0.01278771487600239
$ws.Cells.Item(12, 16).Value = 0.01278771487600239
$ws.Cells.Item(12, 17).Value = 12.71791743768
$ws.Cells.Item(12, 18).Value = 114.46125693912
$ws.Cells.Item(12, 19).Value = 0.005232029223288252
$ws.Cells.Item(12, 20).Value = 0.005232029223288251
$ws.Cells.Item(13, 7).Value = 150.629115
$ws.Cells.Item(13, 8).Value = 451.887345
$ws.Cells.Item(13, 9).Value = 0.4091449703110563
$ws.Cells.Item(13, 10).Value = 0.4091449703110563
$ws.Cells.Item(13, 13).Value = 6.346903333333334
$ws.Cells.Item(13, 14).Value = 19.04071
$ws.Cells.Item(13, 15).Value = 0.9612752294416316
$ws.Cells.Item(13, 16).Value = 0.9612752294416316
$ws.Cells.Item(13, 17).Value = 956.0284320905499
$ws.Cells.Item(13, 18).Value = 8604.255888814951
$ws.Cells.Item(13, 19).Value = 0.3933009252106502
$ws.Cells.Item(13, 20).Value = 0.3933009252106502
$ws.Cells.Item(14, 7).Value = 54.80491966666667
$ws.Cells.Item(14, 8).Value = 164.414759
$ws.Cells.Item(14, 9).Value = 0.1488633670185088
$ws.Cells.Item(14, 10).Value = 0.1488633670185088
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.06778666666666666
$ws.Cells.Item(14, 14).Value = 0.20336
$ws.Cells.Item(14, 15).Value = 0.01026668284214455
$ws.Cells.Item(14, 16).Value = 0.01026668284214455
$ws.Cells.Item(14, 17).Value = 3.715042821137778
$ws.Cells.Item(14, 18).Value = 33.43538539024
$ws.Cells.Item(14, 19).Value = 0.001528332975992791
$ws.Cells.Item(14, 20).Value = 0.001528332975992791
$ws.Cells.Item(15, 7).Value = 54.80491966666667
$ws.Cells.Item(15, 8).Value = 164.414759
$ws.Cells.Item(15, 9).Value = 0.1488633670185088
$ws.Cells.Item(15, 10).Value = 0.1488633670185088
$ws.Cells.Item(15, 15).Value = 0.01567037284022157
$ws.Cells.Item(15, 16).Value = 0.01567037284022157
$ws.Cells.Item(15, 17).Value = 5.670391013311666
$ws.Cells.Item(15, 18).Value = 51.033519119805
$ws.Cells.Item(15, 19).Value = 0.002332744463430775
$ws.Cells.Item(15, 20).Value = 0.002332744463430775
$ws.Cells.Item(16, 7).Value = 54.80491966666667
$ws.Cells.Item(16, 8).Value = 164.414759
$ws.Cells.Item(16, 9).Value = 0.1488633670185088
$ws.Cells.Item(16, 10).Value = 0.1488633670185088
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.08443199999999999
$ws.Cells.Item(16, 14).Value = 0.253296
$ws.Cells.Item(16, 15).Value = 0.01278771487600239
$ws.Cells.Item(16, 16).Value = 0.01278771487600239
$ws.Cells.Item(16, 17).Value = 4.627288977296
$ws.Cells.Item(16, 18).Value = 41.645600795664
$ws.Cells.Item(16, 19).Value = 0.001903622292914388
$ws.Cells.Item(16, 20).Value = 0.001903622292914388
$ws.Cells.Item(17, 7).Value = 54.80491966666667
$ws.Cells.Item(17, 8).Value = 164.414759
$ws.Cells.Item(17, 9).Value = 0.1488633670185088
$ws.Cells.Item(17, 10).Value = 0.1488633670185088
$ws.Cells.Item(17, 13).Value = 6.346903333333334
$ws.Cells.Item(17, 14).Value = 19.04071
$ws.Cells.Item(17, 15).Value = 0.9612752294416316
$ws.Cells.Item(17, 16).Value = 0.9612752294416316
$ws.Cells.Item(17, 17).Value = 347.8415273154322
$ws.Cells.Item(17, 18).Value = 3130.57374583889
$ws.Cells.Item(17, 19).Value = 0.1430986672861708
$ws.Cells.Item(17, 20).Value = 0.1430986672861708